$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Configuration")

# TemplateDefinition.name is now TemplateDefinition.Id:
# the row that used to hold the "Name" variable (Letter-Template) is now
# labeled "Id", and gets a real description instead of being blank.
$ws.Range("A4").Value = "Id"
$ws.Range("C4").Value = "Find a unique name shortly describing the functionality of this template, e. g. 'Employee contract'. You may refer this definition Excel file by this name."

# The "Description" row also gets a real description instead of being blank.
$ws.Range("C5").Value = "Only for describing the purpose of this template definition for the users."

# TemplateDefinition.id isn't a generated value anymore, so the old "Id" row
# (random generated id + "Please do not modify this value." warning) is gone.
$ws.Rows(7).Delete()

# Column C needs to be widened to fit the new, longer description text.
$ws.Columns("C").ColumnWidth = 115.4986979166667

# Focus moves from the Configuration sheet back to the Variables sheet.
$wsVariables = $wb.Worksheets.Item("Variables")
$wsVariables.Activate() | Out-Null
$wsVariables.Range("A1:H1").Select() | Out-Null
